$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Currency number format used by the "Valor" column (style index 1 in styles.xml)
$euroFormat = '#,##0.00\ "€"'

# Apply the consistent euro number format to cells that previously used a
# different/no format (C1 header gets the same numFmt as the rest of the
# column; C13/C14/C16 pick up the euro format; C15 drops the red-negative
# variant in favour of the plain one).
$ws.Range("C1").NumberFormat = $euroFormat
$ws.Range("C13").NumberFormat = $euroFormat
$ws.Range("C14").NumberFormat = $euroFormat
$ws.Range("C15").NumberFormat = $euroFormat
$ws.Range("C16").NumberFormat = $euroFormat

# Price correction for "header 8 pinos"
$ws.Range("C9").Value = 0.9

# New purchased item: "interruptor"
$ws.Range("A17").Value = "interruptor"
$ws.Range("C17").Value = 0.65
$ws.Range("C17").NumberFormat = $euroFormat

# Move the active selection like in the edited workbook
$ws.Range("E19").Select()
